# Applies the "Fixed update to excel issue" change:
#  1. Rename the "Requested quantity" header on the Weekly Quantity and
#     Monthly Trend sheets to Weekly_PO_Qty / Monthly_PO_Qty respectively.
#  2. Add a new "PO Forecast" worksheet (after "Monthly Trend") containing a
#     forecast table (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1. Rename the header cells -------------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" -------------------
$newSheet = $wb.Worksheets.Add($null, $wsMonthly)
$newSheet.Name = "PO Forecast"

# Reuse the existing header style (bold + border, same as the other sheets'
# header row) by copying formats from the Weekly Quantity header row.
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Reuse the existing date-cell style (numFmt 165) for column A data rows.
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A11").PasteSpecial(-4122)

$newSheet.Range("A2").Value  = 45578.99999999999
$newSheet.Range("B2").Value  = 456
$newSheet.Range("C2").Value  = 456.0017794695483
$newSheet.Range("D2").Value  = 456.0017795501641

$newSheet.Range("A3").Value  = 45592.99999999999
$newSheet.Range("B3").Value  = 24
$newSheet.Range("C3").Value  = 24.00197164134324
$newSheet.Range("D3").Value  = 24.00197172197931

$newSheet.Range("A4").Value  = 45599.99999999999
$newSheet.Range("B4").Value  = 0
$newSheet.Range("C4").Value  = -191.9979331708318
$newSheet.Range("D4").Value  = -191.9979313255052

$newSheet.Range("A5").Value  = 45606.99999999999
$newSheet.Range("B5").Value  = 0
$newSheet.Range("C5").Value  = -407.9978396821843
$newSheet.Range("D5").Value  = -407.997832929343

$newSheet.Range("A6").Value  = 45613.99999999999
$newSheet.Range("B6").Value  = 0
$newSheet.Range("C6").Value  = -623.997747782146
$newSheet.Range("D6").Value  = -623.9977332328905

$newSheet.Range("A7").Value  = 45620.99999999999
$newSheet.Range("B7").Value  = 0
$newSheet.Range("C7").Value  = -839.9976564548732
$newSheet.Range("D7").Value  = -839.9976332263219

$newSheet.Range("A8").Value  = 45627.99999999999
$newSheet.Range("B8").Value  = 0
$newSheet.Range("C8").Value  = -1055.997565849447
$newSheet.Range("D8").Value  = -1055.997531458588

$newSheet.Range("A9").Value  = 45634.99999999999
$newSheet.Range("B9").Value  = 0
$newSheet.Range("C9").Value  = -1271.997476724994
$newSheet.Range("D9").Value  = -1271.997429045829

$newSheet.Range("A10").Value = 45641.99999999999
$newSheet.Range("B10").Value = 0
$newSheet.Range("C10").Value = -1487.997387900627
$newSheet.Range("D10").Value = -1487.997326543177

$newSheet.Range("A11").Value = 45648.99999999999
$newSheet.Range("B11").Value = 0
$newSheet.Range("C11").Value = -1703.997298532983
$newSheet.Range("D11").Value = -1703.997223481789
